$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18; existing rows 18-30 shift down to 19-31
$ws.Rows("18:18").Insert()

# Populate the new row 18 with the new weekly data point
$ws.Range("A18").Value = 10
$ws.Range("B18").Value = "Vega Modelo de Temuco"
$ws.Range("C18").Value = "La Araucanía"
$ws.Range("D18").Value = 45240
$ws.Range("D18").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E18").Value = 9
$ws.Range("F18").Value = "Fruta"
$ws.Range("G18").Value = 100104
$ws.Range("H18").Value = "Frutos de pepita"
$ws.Range("I18").Value = 100104004
$ws.Range("J18").Value = "Níspero"
$ws.Range("K18").Value = "Californiana(o)"
$ws.Range("L18").Value = "Primera"
$ws.Range("M18").Value = 80
$ws.Range("N18").Value = 24000
$ws.Range("O18").Value = 24000
$ws.Range("P18").Value = 24000
$ws.Range("Q18").Value = "$/bandeja 5 kilos"
$ws.Range("R18").Value = "Provincia de Quillota"
$ws.Range("S18").Value = 4800
$ws.Range("T18").Value = 5
